{"js": "// Update the date line and every two-digit \u00f7 one-digit practice\n// answer cell in the table to the new day's values.\nconst replacements = [\n  [\"2025-08-16 Saturday\", \"2025-08-17 Sunday\"],\n  [\"28\u00f72=14, 0\", \"30\u00f77=4, 2\"],\n  [\"77\u00f76=12, 5\", \"42\u00f77=6, 0\"],\n  [\"14\u00f75=2, 4\", \"45\u00f79=5, 0\"],\n  [\"83\u00f74=20, 3\", \"84\u00f73=28, 0\"],\n  [\"72\u00f76=12, 0\", \"82\u00f78=10, 2\"],\n  [\"35\u00f76=5, 5\", \"35\u00f74=8, 3\"],\n  [\"12\u00f73=4, 0\", \"94\u00f76=15, 4\"],\n  [\"26\u00f74=6, 2\", \"25\u00f73=8, 1\"],\n  [\"11\u00f79=1, 2\", \"36\u00f78=4, 4\"],\n  [\"81\u00f75=16, 1\", \"19\u00f74=4, 3\"],\n  [\"43\u00f73=14, 1\", \"65\u00f76=10, 5\"],\n  [\"74\u00f76=12, 2\", \"21\u00f72=10, 1\"],\n  [\"89\u00f78=11, 1\", \"66\u00f74=16, 2\"],\n  [\"57\u00f77=8, 1\", \"98\u00f75=19, 3\"],\n  [\"62\u00f79=6, 8\", \"76\u00f76=12, 4\"],\n  [\"28\u00f74=7, 0\", \"81\u00f79=9, 0\"],\n  [\"11\u00f73=3, 2\", \"24\u00f78=3, 0\"],\n  [\"85\u00f78=10, 5\", \"65\u00f74=16, 1\"],\n  [\"68\u00f79=7, 5\", \"13\u00f79=1, 4\"],\n  [\"24\u00f75=4, 4\", \"43\u00f78=5, 3\"],\n  [\"65\u00f77=9, 2\", \"28\u00f73=9, 1\"],\n  [\"71\u00f76=11, 5\", \"74\u00f73=24, 2\"],\n  [\"73\u00f75=14, 3\", \"57\u00f75=11, 2\"],\n  [\"66\u00f73=22, 0\", \"41\u00f77=5, 6\"],\n  [\"44\u00f73=14, 2\", \"13\u00f75=2, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every two-digit \u00f7 one-digit practice\n# answer cell in the table to the new day's values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-08-16 Saturday\", \"2025-08-17 Sunday\"),\n    @(\"28\u00f72=14, 0\", \"30\u00f77=4, 2\"),\n    @(\"77\u00f76=12, 5\", \"42\u00f77=6, 0\"),\n    @(\"14\u00f75=2, 4\", \"45\u00f79=5, 0\"),\n    @(\"83\u00f74=20, 3\", \"84\u00f73=28, 0\"),\n    @(\"72\u00f76=12, 0\", \"82\u00f78=10, 2\"),\n    @(\"35\u00f76=5, 5\", \"35\u00f74=8, 3\"),\n    @(\"12\u00f73=4, 0\", \"94\u00f76=15, 4\"),\n    @(\"26\u00f74=6, 2\", \"25\u00f73=8, 1\"),\n    @(\"11\u00f79=1, 2\", \"36\u00f78=4, 4\"),\n    @(\"81\u00f75=16, 1\", \"19\u00f74=4, 3\"),\n    @(\"43\u00f73=14, 1\", \"65\u00f76=10, 5\"),\n    @(\"74\u00f76=12, 2\", \"21\u00f72=10, 1\"),\n    @(\"89\u00f78=11, 1\", \"66\u00f74=16, 2\"),\n    @(\"57\u00f77=8, 1\", \"98\u00f75=19, 3\"),\n    @(\"62\u00f79=6, 8\", \"76\u00f76=12, 4\"),\n    @(\"28\u00f74=7, 0\", \"81\u00f79=9, 0\"),\n    @(\"11\u00f73=3, 2\", \"24\u00f78=3, 0\"),\n    @(\"85\u00f78=10, 5\", \"65\u00f74=16, 1\"),\n    @(\"68\u00f79=7, 5\", \"13\u00f79=1, 4\"),\n    @(\"24\u00f75=4, 4\", \"43\u00f78=5, 3\"),\n    @(\"65\u00f77=9, 2\", \"28\u00f73=9, 1\"),\n    @(\"71\u00f76=11, 5\", \"74\u00f73=24, 2\"),\n    @(\"73\u00f75=14, 3\", \"57\u00f75=11, 2\"),\n    @(\"66\u00f73=22, 0\", \"41\u00f77=5, 6\"),\n    @(\"44\u00f73=14, 2\", \"13\u00f75=2, 3\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
